$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6.. down by one.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new record's data.
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Femacal de La Calera"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44503
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 100112026
$ws.Range("G6").Value = "Haba"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 73
$ws.Range("K6").Value = 7500
$ws.Range("L6").Value = 8000
$ws.Range("M6").Value = 7740
$ws.Range("N6").Value = '$/saco 25 kilos'
$ws.Range("O6").Value = "Provincia de Quillota"
$ws.Range("P6").Value = 310
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
